# fix(excel): correction d'une donnée manquante dans les exports
# Adds the two missing "orientation" rows (CIAS / Autre orientation) to the
# "Répartition des orientations" block of the structure-stats export, just
# below "Orientation vers Organisme agrée" (row 101) and above the
# "3. TOTAL DES INTERACTIONS" section header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102 already exists but was left blank -- insert one more row below it
# so the previously-blank row keeps its place and a brand new row 103 shows
# up for the second missing label. Everything from the old row 103 onward
# (the "3. TOTAL DES INTERACTIONS..." block) shifts down by one row.
$ws.Rows("103:103").Insert()

# Fill in the two newly-available rows with the missing orientation labels.
$ws.Range("B102").Value = "Orientation vers CIAS"
$ws.Range("B103").Value = "Autre orientation"

# Match the row height used by all the other sibling rows in this block.
$ws.Rows("102:102").RowHeight = 16
$ws.Rows("103:103").RowHeight = 16

# Restore the view: scrolled up a bit and the selection left on B101.
$ws.Range("B101").Select()
